$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — update "想去人数" (F column) counts.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 6272
$wsExhibit.Range("F4").Value = 557
$wsExhibit.Range("F5").Value = 121
$wsExhibit.Range("F6").Value = 26
$wsExhibit.Range("F8").Value = 342
$wsExhibit.Range("F9").Value = 1331
$wsExhibit.Range("F10").Value = 92

# Sheet "全部类型" (All types) — same underlying rows, different row numbers.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 6272
$wsAll.Range("F4").Value = 557
$wsAll.Range("F5").Value = 121
$wsAll.Range("F6").Value = 26
$wsAll.Range("F8").Value = 342
$wsAll.Range("F13").Value = 1331
$wsAll.Range("F14").Value = 92
